$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.232.20'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.779.29'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '594.96'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '167.46'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").Value = '3.777.15'
$ws.Range("E7").Value = '  +1.00%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '6.40'
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '0.0000259'
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("D14").Value = '36.12'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '4.409.30'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '3.749.48'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '68.220.83'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '17.87'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").Value = '10.82'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '464.56'
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").Value = '0.698'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  +9.94%  '
$ws.Range("D25").Value = '83.91'
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("D27").Value = '11.84'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("D28").Value = '10.15'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").Value = '29.93'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").Value = '2.16'
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("D34").Value = '9.14'
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '3.734.05'
$ws.Range("E36").Value = '  +1.22%  '
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").Value = '5.78'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D44").Value = '44.42'
$ws.Range("E44").Value = '  +16.44%  '
$ws.Range("E45").Value = '  -2.34%  '
$ws.Range("D46").Value = '46.97'
$ws.Range("E46").Value = '  +3.42%  '
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").Value = '8.41'
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("D49").Value = '145.83'
$ws.Range("E49").Value = '  +1.61%  '
$ws.Range("D50").Value = '391.58'
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = '2.783.83'
$ws.Range("E51").Value = '  +3.93%  '
